$wb = $excel.ActiveWorkbook

# Sheet ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 33333854
$ws.Range("I6").Value = 40000424
$ws.Range("K6").Value = 120001272
$ws.Range("M6").Value = -120001160

# Sheet ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1821.75
$ws.Range("I19").Value = 1214.6364
$ws.Range("J19").Value = 2563.7778
$ws.Range("K19").Value = 1214.6364
$ws.Range("L19").Value = 2563.7778
$ws.Range("M19").Value = -1039.6364
$ws.Range("N19").Value = -2913.7778

# Sheet ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 11589.8
$ws.Range("J64").Value = 9816.5
$ws.Range("L64").Value = 9816.5
$ws.Range("N64").Value = -10312.5

# Sheet ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 11589.8
$ws.Range("J67").Value = 9816.5
$ws.Range("L67").Value = 9816.5
$ws.Range("N67").Value = -11532.5

# Sheet ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 568.1177
$ws.Range("I92").Value = 413.91666
$ws.Range("J92").Value = 938.2
$ws.Range("K92").Value = 413.91666
$ws.Range("L92").Value = 938.2
$ws.Range("M92").Value = 834.08334
$ws.Range("N92").Value = -3434.2

# Sheet ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1148.6
$ws.Range("I98").Value = 935.75
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 935.75
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 562.25
$ws.Range("N98").Value = -4996

# Sheet ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 543.1177
$ws.Range("I99").Value = 565.8333
$ws.Range("J99").Value = 488.6
$ws.Range("K99").Value = 1697.4999
$ws.Range("L99").Value = 1465.8
$ws.Range("M99").Value = -199.4999
$ws.Range("N99").Value = -4461.8

# Sheet ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2209.842
$ws.Range("I100").Value = 845.7692
$ws.Range("J100").Value = 5165.3335
$ws.Range("K100").Value = 845.7692
$ws.Range("L100").Value = 5165.3335
$ws.Range("M100").Value = -304.7692
$ws.Range("N100").Value = -6247.3335

# Sheet ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1148.6
$ws.Range("I122").Value = 935.75
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2807.25
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -357.25
$ws.Range("N122").Value = -10900

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3931.488
$ws.Range("I32").Value = 2202.7124
$ws.Range("K32").Value = 2202.7124
$ws.Range("M32").Value = -1915.7124

# Sheet ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 8349.714
$ws.Range("I88").Value = 25796.5
$ws.Range("J88").Value = 1371
$ws.Range("K88").Value = 25796.5
$ws.Range("L88").Value = 1371
$ws.Range("M88").Value = -25390.5
$ws.Range("N88").Value = -2183

# Sheet ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 8349.714
$ws.Range("I91").Value = 25796.5
$ws.Range("J91").Value = 1371
$ws.Range("K91").Value = 25796.5
$ws.Range("L91").Value = 1371
$ws.Range("M91").Value = -24392.5
$ws.Range("N91").Value = -4179

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2712.3428
$ws.Range("I122").Value = 2232.0588
$ws.Range("K122").Value = 6696.176399999999
$ws.Range("M122").Value = -4246.176399999999

# Sheet BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2185.9524
$ws.Range("I20").Value = 2094.0715
$ws.Range("J20").Value = 2369.7144
$ws.Range("K20").Value = 2094.0715
$ws.Range("L20").Value = 2369.7144
$ws.Range("M20").Value = -1847.0715
$ws.Range("N20").Value = -2863.7144

# Sheet BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17292.15
$ws.Range("I86").Value = 8138.706
$ws.Range("J86").Value = 69161.664
$ws.Range("K86").Value = 8138.706
$ws.Range("L86").Value = 69161.664
$ws.Range("M86").Value = -7015.706
$ws.Range("N86").Value = -71407.664

# Sheet BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 17292.15
$ws.Range("I89").Value = 8138.706
$ws.Range("J89").Value = 69161.664
$ws.Range("K89").Value = 40693.53
$ws.Range("L89").Value = 345808.32
$ws.Range("M89").Value = -35077.53
$ws.Range("N89").Value = -357040.32

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3835.2727
$ws.Range("I134").Value = 3470.3333
$ws.Range("K134").Value = 10410.9999
$ws.Range("M134").Value = -7875.999899999999

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3293.4285
$ws.Range("I31").Value = 1978.64
$ws.Range("J31").Value = 14250
$ws.Range("K31").Value = 1978.64
$ws.Range("L31").Value = 14250
$ws.Range("M31").Value = -1683.64
$ws.Range("N31").Value = -14840

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3293.4285
$ws.Range("I34").Value = 1978.64
$ws.Range("J34").Value = 14250
$ws.Range("K34").Value = 1978.64
$ws.Range("L34").Value = 14250
$ws.Range("M34").Value = -1776.64
$ws.Range("N34").Value = -14654

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3388.5715
$ws.Range("I58").Value = 2524.1667
$ws.Range("J58").Value = 4036.875
$ws.Range("K58").Value = 2524.1667
$ws.Range("L58").Value = 4036.875
$ws.Range("M58").Value = -2321.1667
$ws.Range("N58").Value = -4442.875

# Sheet CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4663
$ws.Range("I86").Value = 4216.1665
$ws.Range("J86").Value = 5199.2
$ws.Range("K86").Value = 4216.1665
$ws.Range("L86").Value = 5199.2
$ws.Range("M86").Value = -3093.1665
$ws.Range("N86").Value = -7445.2

# Sheet CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4663
$ws.Range("I89").Value = 4216.1665
$ws.Range("J89").Value = 5199.2
$ws.Range("K89").Value = 21080.8325
$ws.Range("L89").Value = 25996
$ws.Range("M89").Value = -15464.8325
$ws.Range("N89").Value = -37228

# Sheet CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 8939.375
$ws.Range("I105").Value = 1929.75
$ws.Range("K105").Value = 1929.75
$ws.Range("M105").Value = -182.75

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1394.4231
$ws.Range("I122").Value = 1431.9375
$ws.Range("K122").Value = 4295.8125
$ws.Range("M122").Value = -1845.8125

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3388.5715
$ws.Range("I136").Value = 2524.1667
$ws.Range("J136").Value = 4036.875
$ws.Range("K136").Value = 7572.500100000001
$ws.Range("L136").Value = 12110.625
$ws.Range("M136").Value = -5022.500100000001
$ws.Range("N136").Value = -17210.625

# Sheet CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1010.2
$ws.Range("I23").Value = 101
$ws.Range("J23").Value = 1237.5
$ws.Range("K23").Value = 303
$ws.Range("L23").Value = 3712.5
$ws.Range("M23").Value = -68
$ws.Range("N23").Value = -4182.5

# Sheet CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 223.52942
$ws.Range("I34").Value = 216.875
$ws.Range("K34").Value = 650.625
$ws.Range("M34").Value = -566.625

# Sheet CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 126.09091
$ws.Range("I38").Value = 30
$ws.Range("J38").Value = 181
$ws.Range("K38").Value = 90
$ws.Range("L38").Value = 543
$ws.Range("M38").Value = 257
$ws.Range("N38").Value = -1237

# Sheet CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 779
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

# Sheet CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2436
$ws.Range("I113").Value = 1767.5
$ws.Range("K113").Value = 5302.5
$ws.Range("M113").Value = -3132.5

# Sheet GSM row 21
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 31263750

# Sheet GSM row 30
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 31263750

# Sheet GSM row 31
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Sheet GSM row 37
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4861.3076
$ws.Range("I102").Value = 4699.6665
$ws.Range("J102").Value = 4999.857
$ws.Range("K102").Value = 4699.6665
$ws.Range("L102").Value = 4999.857
$ws.Range("M102").Value = -3077.6665
$ws.Range("N102").Value = -8243.857

# Sheet GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3094.739
$ws.Range("I113").Value = 1594.4615
$ws.Range("K113").Value = 1594.4615
$ws.Range("M113").Value = 575.5385000000001

# Sheet GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6128.1055
$ws.Range("I126").Value = 7028.125
$ws.Range("J126").Value = 5473.5454
$ws.Range("K126").Value = 21084.375
$ws.Range("L126").Value = 16420.6362
$ws.Range("M126").Value = -18614.375
$ws.Range("N126").Value = -21360.6362

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2571.64
$ws.Range("I132").Value = 2240.875
$ws.Range("K132").Value = 6722.625
$ws.Range("M132").Value = -4192.625

# Sheet LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3839.7222
$ws.Range("I7").Value = 3746.4443
$ws.Range("J7").Value = 3933
$ws.Range("K7").Value = 3746.4443
$ws.Range("L7").Value = 3933
$ws.Range("M7").Value = -3634.4443
$ws.Range("N7").Value = -4157

# Sheet LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3916.9524
$ws.Range("I40").Value = 3793.25
$ws.Range("K40").Value = 3793.25
$ws.Range("M40").Value = -3657.25

# Sheet LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3195.1904
$ws.Range("I100").Value = 2793.2144
$ws.Range("K100").Value = 2793.2144
$ws.Range("M100").Value = -2252.2144

# Sheet LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3839.7222
$ws.Range("I126").Value = 3746.4443
$ws.Range("J126").Value = 3933
$ws.Range("K126").Value = 11239.3329
$ws.Range("L126").Value = 11799
$ws.Range("M126").Value = -8769.332900000001
$ws.Range("N126").Value = -16739

# Sheet WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7785.7144
$ws.Range("J62").Value = 7785.7144
$ws.Range("L62").Value = 7785.7144
$ws.Range("N62").Value = -9033.714400000001

# Sheet WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7785.7144
$ws.Range("J65").Value = 7785.7144
$ws.Range("L65").Value = 38928.572
$ws.Range("N65").Value = -45168.572

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6362.25
$ws.Range("I132").Value = 6603.8823
$ws.Range("J132").Value = 4993
$ws.Range("K132").Value = 19811.6469
$ws.Range("L132").Value = 14979
$ws.Range("M132").Value = -17281.6469
$ws.Range("N132").Value = -20039
